# Update the account-statement table: remove the previous employee/period
# records and replace them with the new ones, grouped by worker and sorted
# by period (1903, 1902, 1901), and refresh the "Salario Basico" value for
# JUAN DARIO LOMBANA HERRERA to 781242.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$docJuan   = "1051818784"
$nameJuan  = "JUAN DARIO LOMBANA HERRERA"
$docDavid  = "84458883"
$nameDavid = "DAVID GUILLERMO RONDON VISBAL"

# Row, DocNumber, Name, Period, ValorMora, SalarioBasico
$rows = @(
    @(16, $docJuan,  $nameJuan,  "1903", 56667, 781242),
    @(17, $docJuan,  $nameJuan,  "1902", 68000, 781242),
    @(18, $docJuan,  $nameJuan,  "1901", 68000, 781242),
    @(19, $docDavid, $nameDavid, "1903", 56667, 1700000),
    @(20, $docDavid, $nameDavid, "1902", 68000, 1700000),
    @(21, $docDavid, $nameDavid, "1901", 68000, 1700000)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("C$rowNum").Value = $r[1]
    $ws.Range("D$rowNum").Value = $r[2]
    $ws.Range("E$rowNum").Value = $r[3]
    $ws.Range("F$rowNum").Value = $r[4]
    $ws.Range("G$rowNum").Value = $r[5]
}
